# Restructure "Fields" sheet: insert two new rows into the model_runs
# section (run_name + model_name), pushing all subsequent rows down by 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fields")

# Insert two blank rows right before the current row 9
# (old row 9 = model_snapshot_date), shifting rows 9-25 down to 11-27.
$ws.Range("A9:A10").EntireRow.Insert()

# --- New row 9: run_name -------------------------------------------------
$ws.Cells.Item(9, 1).Value = "model_runs"
$ws.Cells.Item(9, 2).Value = "run_name"
$ws.Cells.Item(9, 3).Value = "A brief description of the model run used to generate legends in plots"
$ws.Cells.Item(9, 4).Value = "any text field (with under 500 characters)"

$ws.Rows.Item(9).RowHeight = 50
$ws.Cells.Item(9, 3).WrapText = $true
$ws.Cells.Item(9, 4).WrapText = $true

# --- New row 10: model_name ----------------------------------------------
$ws.Cells.Item(10, 1).Value = "model_runs"
$ws.Cells.Item(10, 2).Value = "model_name"
$ws.Cells.Item(10, 3).Value = "The name of each model"
$ws.Cells.Item(10, 4).Value = "any text field (with under 500 characters)"

$ws.Rows.Item(10).RowHeight = 50
$ws.Cells.Item(10, 3).WrapText = $true
$ws.Cells.Item(10, 4).WrapText = $true

# Move the active selection to D10, matching the authored edit.
[void]$ws.Range("D10").Select()
